$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $cellRef, $value, $styleSourceRef)
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value2 = $value
    $ws.Range($cellRef).Style = $ws.Range($styleSourceRef).Style
}

Set-CellText $ws 'D2' '69.718.91' 'D31'
Set-CellText $ws 'D3' '3.927.76' 'D31'
Set-CellText $ws 'E3' '  +0.26%  ' 'E31'
Set-CellText $ws 'D4' '1.00' 'D31'
Set-CellText $ws 'E4' '  +0.02%  ' 'E31'
Set-CellText $ws 'D5' '528.11' 'D31'
Set-CellText $ws 'E5' '  +8.75%  ' 'E31'
Set-CellText $ws 'D6' '144.48' 'D31'
Set-CellText $ws 'E6' '  -1.10%  ' 'E31'
Set-CellText $ws 'E7' '  -1.17%  ' 'E31'
Set-CellText $ws 'E8' '  +0.03%  ' 'E31'
Set-CellText $ws 'D9' '0.724' 'D31'
Set-CellText $ws 'E9' '  -1.46%  ' 'E31'
Set-CellText $ws 'D10' '0.172' 'D31'
Set-CellText $ws 'E10' '  +3.28%  ' 'E31'
Set-CellText $ws 'D11' '0.0000334' 'D31'
Set-CellText $ws 'E11' '  -3.09%  ' 'E31'
Set-CellText $ws 'D12' '42.35' 'D31'
Set-CellText $ws 'E12' '  -2.24%  ' 'E31'
Set-CellText $ws 'D13' '4.557.09' 'D31'
Set-CellText $ws 'E13' '  +0.46%  ' 'E31'
Set-CellText $ws 'D14' '10.32' 'D31'
Set-CellText $ws 'E14' '  -4.73%  ' 'E31'
Set-CellText $ws 'D15' '3.921.84' 'D31'
Set-CellText $ws 'E15' '  -0.09%  ' 'E31'
Set-CellText $ws 'B16' 'TRON' 'B31'
Set-CellText $ws 'C16' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx' 'C31'
Set-CellText $ws 'D16' '0.136' 'D31'
Set-CellText $ws 'E16' '  -0.31%  ' 'E31'
Set-CellText $ws 'B17' 'Uniswap' 'B31'
Set-CellText $ws 'C17' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni' 'C31'
Set-CellText $ws 'D17' '13.93' 'D31'
Set-CellText $ws 'E17' '  -2.72%  ' 'E31'
Set-CellText $ws 'D18' '1.22' 'D31'
Set-CellText $ws 'E18' '  +6.79%  ' 'E31'
Set-CellText $ws 'E19' '  -1.78%  ' 'E31'
Set-CellText $ws 'D20' '69.592.09' 'D31'
Set-CellText $ws 'E20' '  +1.74%  ' 'E31'
Set-CellText $ws 'D21' '432.48' 'D31'
Set-CellText $ws 'E22' '  -3.93%  ' 'E31'
Set-CellText $ws 'D23' '14.41' 'D31'
Set-CellText $ws 'E23' '  -4.33%  ' 'E31'
Set-CellText $ws 'D24' '4.13' 'D31'
Set-CellText $ws 'E24' '  +14.52%  ' 'E31'
Set-CellText $ws 'D25' '87.87' 'D31'
Set-CellText $ws 'E25' '  -0.40%  ' 'E31'
Set-CellText $ws 'D26' '11.68' 'D31'
Set-CellText $ws 'E26' '  +3.83%  ' 'E31'
Set-CellText $ws 'D27' '10.72' 'D31'
Set-CellText $ws 'E27' '  -4.36%  ' 'E31'
Set-CellText $ws 'D28' '36.47' 'D31'
Set-CellText $ws 'E28' '  -4.13%  ' 'E31'
Set-CellText $ws 'D29' '696.90' 'D31'
Set-CellText $ws 'E29' '  -2.30%  ' 'E31'
Set-CellText $ws 'D30' '13.20' 'D31'
Set-CellText $ws 'E30' '  -4.29%  ' 'E31'
Set-CellText $ws 'E32' '  -3.53%  ' 'E31'
Set-CellText $ws 'D33' '68.97' 'D31'
Set-CellText $ws 'E33' '  +13.38%  ' 'E31'
Set-CellText $ws 'D34' '0.457' 'D31'
Set-CellText $ws 'E34' '  +16.30%  ' 'E31'
Set-CellText $ws 'D35' '5.99' 'D31'
Set-CellText $ws 'E35' '  -3.74%  ' 'E31'
Set-CellText $ws 'D36' '40.04' 'D31'
Set-CellText $ws 'E36' '  -3.51%  ' 'E31'
Set-CellText $ws 'D37' '0.0₃0842' 'D31'
Set-CellText $ws 'E37' '  -3.84%  ' 'E31'
Set-CellText $ws 'D38' '0.147' 'D31'
Set-CellText $ws 'E38' '  +1.17%  ' 'E31'
Set-CellText $ws 'D39' '1.00' 'D31'
Set-CellText $ws 'E39' '  +0.02%  ' 'E31'
Set-CellText $ws 'E40' '  -0.01%  ' 'E31'
Set-CellText $ws 'D41' '0.0480' 'D31'
Set-CellText $ws 'E41' '  -1.64%  ' 'E31'
Set-CellText $ws 'D42' '3.11' 'D31'
Set-CellText $ws 'E42' '  +4.81%  ' 'E31'
Set-CellText $ws 'D43' '2.78' 'D31'
Set-CellText $ws 'E43' '  -5.11%  ' 'E31'
Set-CellText $ws 'B44' 'ThetaToken' 'B31'
Set-CellText $ws 'C44' 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta' 'C31'
Set-CellText $ws 'D44' '2.96' 'D31'
Set-CellText $ws 'E44' '  -5.26%  ' 'E31'
Set-CellText $ws 'B45' 'Stacks' 'B31'
Set-CellText $ws 'C45' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx' 'C31'
Set-CellText $ws 'D45' '3.17' 'D31'
Set-CellText $ws 'E45' '  +11.97%  ' 'E31'
Set-CellText $ws 'E46' '  -0.48%  ' 'E31'
Set-CellText $ws 'D47' '0.142' 'D31'
Set-CellText $ws 'E47' '  -0.27%  ' 'E31'
Set-CellText $ws 'D48' '0.0₆0352' 'D31'
Set-CellText $ws 'E48' '  +4.39%  ' 'E31'
Set-CellText $ws 'E49' '  -3.54%  ' 'E31'
Set-CellText $ws 'D50' '144.36' 'D31'
Set-CellText $ws 'E50' '  -0.63%  ' 'E31'
Set-CellText $ws 'E51' '  -3.19%  ' 'E31'
